$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a second Train/Valid/Test phase table (rows 10-15) below the existing
# one (rows 1-8), mirroring its layout/formatting, with a new set of
# quarters (3-year-2-month window) and a helper index block in I:L.
# ---------------------------------------------------------------------------

# --- Row 10: header (copy straight from row 1 so formatting matches) ------
$ws.Range("A1").Copy($ws.Range("A10"))
$ws.Range("B1").Copy($ws.Range("B10"))
$ws.Range("C1").Copy($ws.Range("C10"))
$ws.Range("D1").Copy($ws.Range("D10"))
$ws.Range("E1").Copy($ws.Range("E10"))
$ws.Range("F1").Copy($ws.Range("F10"))
$ws.Range("G1").Copy($ws.Range("G10"))

# --- Row 11 -----------------------------------------------------------------
$ws.Range("A2").Copy($ws.Range("A11"))

$ws.Range("B2").Copy($ws.Range("B11"))
$ws.Range("B11").Value = "2016_Q1"

$ws.Range("C2").Copy($ws.Range("C11"))
$ws.Range("C11").Value = "2019_Q3"

$ws.Range("E2").Copy($ws.Range("D11"))
$ws.Range("D11").Value = "2019_Q3"

$ws.Range("D2").Copy($ws.Range("E11"))
$ws.Range("E11").Value = "2020_Q3"

$ws.Range("F2").Copy($ws.Range("F11"))
$ws.Range("F11").Value = "2020_Q3"

$ws.Range("G2").Copy($ws.Range("G11"))
$ws.Range("G11").Value = "2021_Q3"

$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 15
$ws.Range("K11").Value = 19
$ws.Range("L11").Value = 23

# --- Row 12 -----------------------------------------------------------------
$ws.Range("A3").Copy($ws.Range("A12"))

$ws.Range("B3").Copy($ws.Range("B12"))
$ws.Range("B12").Value = "2017_Q1"

$ws.Range("C3").Copy($ws.Range("C12"))
$ws.Range("C12").Value = "2020_Q3"

$ws.Range("D3").Copy($ws.Range("D12"))
$ws.Range("D12").Value = "2020_Q3"

$ws.Range("E3").Copy($ws.Range("E12"))
$ws.Range("E12").Value = "2021_Q3"

$ws.Range("F3").Copy($ws.Range("F12"))
$ws.Range("F12").Value = "2021_Q3"

$ws.Range("G3").Copy($ws.Range("G12"))
$ws.Range("G12").Value = "2022_Q3"

$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 19
$ws.Range("K12").Value = 23
$ws.Range("L12").Value = 27

# --- Row 13 -----------------------------------------------------------------
$ws.Range("A4").Copy($ws.Range("A13"))

$ws.Range("B4").Copy($ws.Range("B13"))
$ws.Range("B13").Value = "2018_Q1"

$ws.Range("C4").Copy($ws.Range("C13"))
$ws.Range("C13").Value = "2021_Q3"

$ws.Range("D4").Copy($ws.Range("D13"))
$ws.Range("D13").Value = "2021_Q3"

$ws.Range("E4").Copy($ws.Range("E13"))
$ws.Range("E13").Value = "2022_Q3"

$ws.Range("F4").Copy($ws.Range("F13"))
$ws.Range("F13").Value = "2022_Q3"

$ws.Range("G4").Copy($ws.Range("G13"))
$ws.Range("G13").Value = "2023_Q3"

$ws.Range("I13").Value = 9
$ws.Range("J13").Value = 23
$ws.Range("K13").Value = 27
$ws.Range("L13").Value = 31

# --- Row 14 -----------------------------------------------------------------
$ws.Range("A5").Copy($ws.Range("A14"))

$ws.Range("B5").Copy($ws.Range("B14"))
$ws.Range("B14").Value = "2019_Q1"

$ws.Range("C5").Copy($ws.Range("C14"))
$ws.Range("C14").Value = "2022_Q3"

$ws.Range("D5").Copy($ws.Range("D14"))
$ws.Range("D14").Value = "2022_Q3"

$ws.Range("E5").Copy($ws.Range("E14"))
$ws.Range("E14").Value = "2023_Q3"

$ws.Range("F5").Copy($ws.Range("F14"))
$ws.Range("F14").Value = "2023_Q3"

$ws.Range("G5").Copy($ws.Range("G14"))
$ws.Range("G14").Value = "2024_Q3"

$ws.Range("I14").Value = 13
$ws.Range("J14").Value = 27
$ws.Range("K14").Value = 31
$ws.Range("L14").Value = 35

# --- Row 15: footer (copy from row 7's merged "기간" summary row) ----------
$ws.Range("A7").Copy($ws.Range("A15"))

$ws.Range("B7").Copy($ws.Range("B15"))
$ws.Range("B15").Value = "3년 2개월"

$ws.Range("C7").Copy($ws.Range("C15"))

$ws.Range("D7").Copy($ws.Range("D15"))
$ws.Range("D15").Value = "1년"

$ws.Range("E7").Copy($ws.Range("E15"))

$ws.Range("F7").Copy($ws.Range("F15"))
$ws.Range("F15").Value = "1년"

$ws.Range("G7").Copy($ws.Range("G15"))

# Merge the footer row's label/value pairs, same as row 7.
$ws.Range("B15:C15").Merge()
$ws.Range("D15:E15").Merge()
$ws.Range("F15:G15").Merge()

$ws.Range("H18").Select()
